$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44365
$ws.Range("J2").Value = 55
$ws.Range("K2").Value = 5000
$ws.Range("L2").Value = 5000
$ws.Range("M2").Value = 5000
$ws.Range("P2").Value = 5000

$ws.Range("D3").Value = 44291
$ws.Range("J3").Value = 35

$ws.Range("D4").Value = 44301
$ws.Range("K4").Value = 3000
$ws.Range("L4").Value = 3000
$ws.Range("M4").Value = 3000
$ws.Range("P4").Value = 3000

$ws.Range("D5").Value = 44497
$ws.Range("J5").Value = 20
$ws.Range("K5").Value = 4000
$ws.Range("L5").Value = 4000
$ws.Range("M5").Value = 4000
$ws.Range("P5").Value = 4000

$ws.Range("D6").Value = 44498
$ws.Range("J6").Value = 40

$ws.Range("D7").Value = 44509
$ws.Range("J7").Value = 20

$ws.Range("D8").Value = 44508

$ws.Range("D9").Value = 44313
$ws.Range("J9").Value = 20

$ws.Range("D10").Value = 44504
$ws.Range("J10").Value = 55

$ws.Range("D11").Value = 44280
$ws.Range("J11").Value = 55

$ws.Range("D12").Value = 44259
$ws.Range("J12").Value = 30

$ws.Range("D13").Value = 44316
$ws.Range("J13").Value = 20

$ws.Range("D14").Value = 44315
$ws.Range("J14").Value = 40

$ws.Range("D15").Value = 44176
$ws.Range("J15").Value = 10
$ws.Range("K15").Value = 4000
$ws.Range("L15").Value = 4000
$ws.Range("M15").Value = 4000
$ws.Range("P15").Value = 4000
